$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the last three data rows (rows 5-7): the updated TPM data only
# has three target clusters for the FAPs -> Tnfsf11/Tnfrsf11b pairing now
# (the MuSCs-sourced rows are gone), so the sheet shrinks to A1:T4.
$ws.Rows("5:7").Delete()

# Row 2: FAPs -> Tnfsf11 -> Tnfrsf11b -> FAPs
$ws.Range("A2").Value = "FAPs"
$ws.Range("B2").Value = "Tnfsf11"
$ws.Range("C2").Value = "Tnfrsf11b"
$ws.Range("D2").Value = "FAPs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.604474
$ws.Range("H2").Value = 4.813422
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 2.214957333333333
$ws.Range("N2").Value = 6.644871999999999
$ws.Range("O2").Value = 0.8812411509483107
$ws.Range("P2").Value = 0.8812411509483107
$ws.Range("Q2").Value = 3.553841452442666
$ws.Range("R2").Value = 31.984573071984
$ws.Range("S2").Value = 0.8812411509483107
$ws.Range("T2").Value = 0.8812411509483107

# Row 3: FAPs -> Tnfsf11 -> Tnfrsf11b -> MuSCs
$ws.Range("A3").Value = "FAPs"
$ws.Range("B3").Value = "Tnfsf11"
$ws.Range("C3").Value = "Tnfrsf11b"
$ws.Range("D3").Value = "MuSCs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.604474
$ws.Range("H3").Value = 4.813422
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.274148
$ws.Range("N3").Value = 0.8224440000000001
$ws.Range("O3").Value = 0.1090723037479928
$ws.Range("P3").Value = 0.1090723037479928
$ws.Range("Q3").Value = 0.439863338152
$ws.Range("R3").Value = 3.958770043368
$ws.Range("S3").Value = 0.1090723037479928
$ws.Range("T3").Value = 0.1090723037479928

# Row 4: FAPs -> Tnfsf11 -> Tnfrsf11b -> Resolving-Mac
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Tnfsf11"
$ws.Range("C4").Value = "Tnfrsf11b"
$ws.Range("D4").Value = "Resolving-Mac"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.604474
$ws.Range("H4").Value = 4.813422
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.02434666666666667
$ws.Range("N4").Value = 0.07303999999999999
$ws.Range("O4").Value = 0.009686545303696538
$ws.Range("P4").Value = 0.009686545303696536
$ws.Range("Q4").Value = 0.03906359365333333
$ws.Range("R4").Value = 0.35157234288
$ws.Range("S4").Value = 0.009686545303696538
$ws.Range("T4").Value = 0.009686545303696536
